$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Cell C9: "Fornecer dados Pessoais " -> "Fornece dados Pessoais "
$ws.Range("C9").Value = "Fornece dados Pessoais "

# 2. Cell A14: update the exception text's bracketed clause
$ws.Range("A14").Value = "Excepção 1               (passo 2)`n[Cliente existente]"

# 3. Row 14 height: 90 -> 60
$ws.Rows.Item(14).RowHeight = 60

# 4. Update selection/active cell to A14
$ws.Range("A14").Select()
